$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3565323912770566
$ws.Range("C2").Value = 0.0870157294438485
$ws.Range("D2").Value = 0.0383407730807761
$ws.Range("F2").Value = 0.4337137820672794
$ws.Range("G2").Value = 0.2828693052985862
$ws.Range("H2").Value = 0.4762879322447766
$ws.Range("K2").Value = 0.3674132846893485
$ws.Range("O2").Value = 1.442789809046218
$ws.Range("B3").Value = 0.3122744505098467
$ws.Range("C3").Value = 0.08211892465156723
$ws.Range("D3").Value = 0.03392169679352008
$ws.Range("F3").Value = 0.4386645749807805
$ws.Range("G3").Value = 0.2883860654683765
$ws.Range("H3").Value = 0.4822841368962756
$ws.Range("K3").Value = 0.3210089477768179
$ws.Range("O3").Value = 1.466838253406372
$ws.Range("B4").Value = 0.2850133965830821
$ws.Range("C4").Value = 0.07911123558817224
$ws.Range("D4").Value = 0.03119365216905123
$ws.Range("F4").Value = 0.4420525848236174
$ws.Range("G4").Value = 0.2920509052225526
$ws.Range("H4").Value = 0.486206100345985
$ws.Range("K4").Value = 0.2923896173732032
$ws.Range("O4").Value = 1.482688745418123
$ws.Range("B5").Value = 0.2738833641419092
$ws.Range("C5").Value = 0.07788543874841025
$ws.Range("D5").Value = 0.03007832426484924
$ws.Range("F5").Value = 0.4435207329425346
$ws.Range("G5").Value = 0.2936140224865085
$ws.Range("H5").Value = 0.4878647608462607
$ws.Range("K5").Value = 0.2806959755789649
$ws.Range("O5").Value = 1.489420455430256
$ws.Range("B6").Value = 0.2720339906698257
$ws.Range("C6").Value = 0.07768189125870606
$ws.Range("D6").Value = 0.02989290795011357
$ws.Range("F6").Value = 0.4437698014660505
$ws.Range("G6").Value = 0.2938777809662696
$ws.Range("H6").Value = 0.4881438302120173
$ws.Range("K6").Value = 0.2787524067297795
$ws.Range("O6").Value = 1.490554702413021
$ws.Range("B7").Value = 0.2848633767373201
$ws.Range("C7").Value = 0.07909470446769262
$ws.Range("D7").Value = 0.03117862505813918
$ws.Range("F7").Value = 0.4420720305403236
$ws.Range("G7").Value = 0.2920717040569372
$ws.Range("H7").Value = 0.4862282249339422
$ws.Range("K7").Value = 0.2922320373204457
$ws.Range("O7").Value = 1.482778428562817
$ws.Range("B8").Value = 0.3412906837952221
$ws.Range("C8").Value = 0.08532759553324354
$ws.Range("D8").Value = 0.0368201715699854
$ws.Range("F8").Value = 0.435348514862401
$ws.Range("G8").Value = 0.2847138114472969
$ws.Range("H8").Value = 0.4783055627035004
$ws.Range("K8").Value = 0.3514399121341683
$ws.Range("O8").Value = 1.450856452893781
$ws.Range("B9").Value = 0.4512273422271562
$ws.Range("C9").Value = 0.09753742352377515
$ws.Range("D9").Value = 0.04776379574951761
$ws.Range("F9").Value = 0.4249290310035434
$ws.Range("G9").Value = 0.2724923605692133
$ws.Range("H9").Value = 0.4646750434847569
$ws.Range("K9").Value = 0.4665067009786696
$ws.Range("O9").Value = 1.39687388064192
$ws.Range("B10").Value = 0.5315269208230404
$ws.Range("C10").Value = 0.1064947644251362
$ws.Range("D10").Value = 0.05572839218545766
$ws.Range("F10").Value = 0.4189632658856084
$ws.Range("G10").Value = 0.2648660387806601
$ws.Range("H10").Value = 0.4558214411721053
$ws.Range("K10").Value = 0.5503761351103265
$ws.Range("O10").Value = 1.362478717769122
$ws.Range("B11").Value = 0.5679483693098462
$ws.Range("C11").Value = 0.1105657407492231
$ws.Range("D11").Value = 0.05933469209526265
$ws.Range("F11").Value = 0.4166169077659561
$ws.Range("G11").Value = 0.2616918757179789
$ws.Range("H11").Value = 0.4520455156788756
$ws.Range("K11").Value = 0.5883778105118722
$ws.Range("O11").Value = 1.347977648447909
$ws.Range("B12").Value = 0.5817240751320867
$ws.Range("C12").Value = 0.1121066635212173
$ws.Range("D12").Value = 0.06069781854259304
$ws.Range("F12").Value = 0.4157813089446236
$ws.Range("G12").Value = 0.2605324636579809
$ws.Range("H12").Value = 0.4506518440511584
$ws.Range("K12").Value = 0.6027455716649683
$ws.Range("O12").Value = 1.342651440820092
$ws.Range("B13").Value = 0.5787579697140188
$ws.Range("C13").Value = 0.1117748299306811
$ws.Range("D13").Value = 0.06040435738478322
$ws.Range("F13").Value = 0.4159589155426673
$ws.Range("G13").Value = 0.2607802680388929
$ws.Range("H13").Value = 0.450950386537464
$ws.Range("K13").Value = 0.5996522394912915
$ws.Range("O13").Value = 1.34379119005375
$ws.Range("B14").Value = 0.5690820367873926
$ws.Range("C14").Value = 0.1106925274929296
$ws.Range("D14").Value = 0.0594468879121024
$ws.Range("F14").Value = 0.4165471016339666
$ws.Range("G14").Value = 0.2615956361273817
$ws.Range("H14").Value = 0.4519301321489166
$ws.Range("K14").Value = 0.5895603141636059
$ws.Range("O14").Value = 1.347536148839069
$ws.Range("B15").Value = 0.5631530986667599
$ws.Range("C15").Value = 0.1100294953164536
$ws.Range("D15").Value = 0.05886008212222293
$ws.Range("F15").Value = 0.4169142758928288
$ws.Range("G15").Value = 0.2621006212627606
$ws.Range("H15").Value = 0.4525349677432899
$ws.Range("K15").Value = 0.5833757410907197
$ws.Range("O15").Value = 1.349851545689404
$ws.Range("B16").Value = 0.529144459883355
$ws.Range("C16").Value = 0.1062286299970197
$ws.Range("D16").Value = 0.05549236638519517
$ws.Range("F16").Value = 0.4191240069540285
$ws.Range("G16").Value = 0.2650794278430268
$ws.Range("H16").Value = 0.456073270620152
$ws.Range("K16").Value = 0.5478895181248333
$ws.Range("O16").Value = 1.363449473509007
$ws.Range("B17").Value = 0.5082531651602551
$ws.Range("C17").Value = 0.1038958682163269
$ws.Range("D17").Value = 0.05342201293727555
$ws.Range("F17").Value = 0.4205737810025312
$ws.Range("G17").Value = 0.2669825094076685
$ws.Range("H17").Value = 0.4583083559745447
$ws.Range("K17").Value = 0.5260805361961332
$ws.Range("O17").Value = 1.372084968351317
$ws.Range("B18").Value = 0.4962270045563173
$ws.Range("C18").Value = 0.1025537786651682
$ws.Range("D18").Value = 0.05222962059944791
$ws.Range("F18").Value = 0.4214422392721104
$ws.Range("G18").Value = 0.2681048784021272
$ws.Range("H18").Value = 0.4596176014871496
$ws.Range("K18").Value = 0.5135224344885501
$ws.Range("O18").Value = 1.37715967422821
$ws.Range("B19").Value = 0.4921534558137068
$ws.Range("C19").Value = 0.1020993145296103
$ws.Range("D19").Value = 0.05182562815552672
$ws.Range("F19").Value = 0.4217422224330107
$ws.Range("G19").Value = 0.2684896579133884
$ws.Range("H19").Value = 0.4600649569960957
$ws.Range("K19").Value = 0.509268075150203
$ws.Range("O19").Value = 1.378896387103417
$ws.Range("B20").Value = 0.5104781247239885
$ws.Range("C20").Value = 0.1041442314878367
$ws.Range("D20").Value = 0.05364256968017855
$ws.Range("F20").Value = 0.4204158699855256
$ws.Range("G20").Value = 0.2667770483391649
$ws.Range("H20").Value = 0.4580679761198354
$ws.Range("K20").Value = 0.5284036110796819
$ws.Range("O20").Value = 1.371154546231381
$ws.Range("B21").Value = 0.5719245424337203
$ws.Range("C21").Value = 0.111010445077909
$ws.Range("D21").Value = 0.05972818842302274
$ws.Range("F21").Value = 0.4163729004416794
$ws.Range("G21").Value = 0.2613549861737212
$ws.Range("H21").Value = 0.4516413748925103
$ws.Range("K21").Value = 0.5925251794301971
$ws.Range("O21").Value = 1.346431682219418
$ws.Range("B22").Value = 0.6119878736706141
$ws.Range("C22").Value = 0.1154939630475553
$ws.Range("D22").Value = 0.0636908767663158
$ws.Range("F22").Value = 0.4140390530394313
$ws.Range("G22").Value = 0.2580595770849499
$ws.Range("H22").Value = 0.4476521619717175
$ws.Range("K22").Value = 0.6342998831738669
$ws.Range("O22").Value = 1.331235942571197
$ws.Range("B23").Value = 0.5906143460797466
$ws.Range("C23").Value = 0.1131014281229881
$ws.Range("D23").Value = 0.06157727921646483
$ws.Range("F23").Value = 0.4152564237571212
$ws.Range("G23").Value = 0.259795640295323
$ws.Range("H23").Value = 0.4497619766379088
$ws.Range("K23").Value = 0.6120163669223473
$ws.Range("O23").Value = 1.339258060495155
$ws.Range("B24").Value = 0.5094722686472153
$ws.Range("C24").Value = 0.1040319494145621
$ws.Range("D24").Value = 0.05354286257006891
$ws.Range("F24").Value = 0.4204871526533793
$ws.Range("G24").Value = 0.2668698492658486
$ws.Range("H24").Value = 0.458176576268837
$ws.Range("K24").Value = 0.5273534106137276
$ws.Range("O24").Value = 1.371574847317817
$ws.Range("B25").Value = 0.4215668316842596
$ws.Range("C25").Value = 0.09423627687593239
$ws.Range("D25").Value = 0.04481633887718317
$ws.Range("F25").Value = 0.4274513525355275
$ws.Range("G25").Value = 0.2755615347103699
$ws.Range("H25").Value = 0.4681585436053268
$ws.Range("K25").Value = 0.4354931591095692
$ws.Range("O25").Value = 1.410553798053201
